$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D17").Value = "Matrix population modeling"
$ws.Range("D5").Value = "Legend of Ambalapuzha; Geometric growth"
$ws.Range("D7").Value = "Stochastic population growth"
$ws.Range("D9").Value = "Basic logistic population growth"
$ws.Range("D11").Value = "Deeper into logistic growth"
$ws.Range("D13").Value = "Life tables"
$ws.Range("D20").Value = "Lemming case study part 1"
$ws.Range("D21").Value = "Lemming case study part 2"

$ws.Range("D29").Select()
